$wb = $excel.ActiveWorkbook

# The original single sheet contains 31 days of data (including one bad
# reading on 2016-12-17 / serial 42721 with an obviously wrong Flow value
# of 0.6). The edit:
#   1. Flags that bad row on the original sheet with a light-cyan fill and
#      renames the sheet to indicate it holds the full 31-count data set.
#   2. Duplicates that sheet, removes the flagged/bad row from the copy,
#      and renames the copy to indicate it holds the corrected 30-count
#      data set. The new sheet becomes the active tab.

$ws31 = $wb.Worksheets.Item(1)

# 1) Highlight the erroneous row (date 2016-12-17, row 21) with a light
#    cyan fill (RGB FFCDFFFF -> OLE/BGR long 16777165) on the original
#    sheet before duplicating it.
$ws31.Range("A21:D21").Interior.Color = 16777165

# 2) Duplicate the sheet right after itself to build the "30ct" version.
$ws31.Copy([System.Reflection.Missing]::Value, $ws31)
$ws30 = $wb.Worksheets.Item(2)

# 3) Remove the bad row from the duplicate - everything below shifts up.
$ws30.Rows.Item(21).Delete()

# 4) Rename both sheets to reflect the row counts.
$ws31.Name = "4_Silt+Clay_Comb_31ct"
$ws30.Name = "4_Silt+Clay_Comb_30ct"

# 5) Give the new sheet the same light-cyan tab color and make it the
#    active / selected sheet (matches activeTab + tabSelected in the diff).
$ws30.Tab.Color = 16777165
[void]$ws30.Range("E32").Select()
